$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Change cell B11 from "R40" to the text "1" (new shared-string entry).
# The leading apostrophe forces Excel to store this numeric-looking
# value as literal text (matching the <c t="s"> cell in the target),
# instead of silently coercing it to the number 1.
$cell = $ws.Range("B11")
$cell.Value = "'1"

$wb.Save()
